$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update row 2 values (new test data); leading apostrophe preserves the
#     original quotePrefix-flavored text style instead of minting a new one.
#     Order matches the shared-string append order of the authored edit. ---
$ws.Range("A2").Value = "'SampleCompany1"
$ws.Range("B2").Value = "'213213216546"
$ws.Range("D2").Value = "'Sam"
$ws.Range("E2").Value = "'Weber"
$ws.Range("F2").Value = "'Sam@sc1.com"
$ws.Range("C2").Value = "'www.samplecompany1.com"
$ws.Range("G2").Value = "'Samuel@123"

# --- Re-apply hyperlink formatting (font/style) to the website + email cells.
#     The sheet already carries hyperlink relationships for C2/F2/G2 from
#     before (left untouched); Hyperlinks.Add is only used momentarily to pick
#     up the built-in "Hyperlink" cell style/font, then the freshly-added
#     (duplicate) link entry is removed again so the link collection itself
#     stays exactly as it was. ---
$ws.Hyperlinks.Add($ws.Range("C2"), "http://www.samplecompany1.com/")
$ws.Hyperlinks.Item($ws.Hyperlinks.Count).Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:Sam@sc1.com")
$ws.Hyperlinks.Item($ws.Hyperlinks.Count).Delete()

# --- Selection moved to C6 ---
$ws.Range("C6").Select()
